# Update Name of Algo
# Applies updated RandomForest-imputed values to the specific cells that
# changed between the before/after versions of the result_data worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D7"  = -7.125500000000003
    "B8"  = 5.739199999999996
    "B10" = 5.0947
    "B12" = 4.779300000000001
    "D14" = -7.8768
    "D15" = -8.139999999999997
    "B18" = 5.651299999999999
    "D18" = -8.471199999999989
    "D20" = -7.547299999999998
    "B25" = 7.069899999999997
    "D29" = -7.2144
    "D30" = -7.412500000000003
    "D31" = -8.443100000000001
    "D35" = -8.244299999999997
    "B37" = 9.517699999999998
    "D40" = -7.802799999999996
    "D44" = -7.325899999999999
    "D50" = -8.213999999999993
    "D54" = -8.0181
    "B55" = 5.421799999999998
    "B68" = 4.582099999999998
    "D68" = -7.164599999999994
    "D76" = -7.092299999999997
    "B77" = 9.595400000000007
    "B78" = 9.468800000000003
    "B79" = 9.003600000000004
    "B80" = 9.619099999999998
    "B81" = 4.735300000000004
    "B82" = 6.698300000000001
    "B84" = 5.096899999999999
    "D87" = -7.799499999999992
    "D88" = -7.352999999999995
    "D92" = -6.123699999999999
    "D96" = -7.487900000000005
    "D98" = -8.502800000000001
    "B101" = 9.621800000000002
    "D101" = -7.496400000000001
    "B102" = 9.653000000000002
    "D102" = -7.470399999999996
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
